$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New seat/name assignments (recalculated report), row 2 = Platz 1 ... row 26 = Platz 25
$names = @(
    "Student100",
    "Student98",
    "Student83",
    "Student108",
    "Student104",
    "Student106",
    "Student91",
    "Student126",
    "Student94",
    "Student102",
    "Student81",
    "Student71",
    "Student87",
    "Student68",
    "Student90",
    "Student89",
    "Student86",
    "Student113",
    "Student79",
    "Student84",
    "Student43",
    "Gernert",
    "Student80",
    "Student99",
    "Student44"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
}
